# Insert a new weekly price record for Zanahoria (Terminal Hortofrutícola Agro
# Chillán) at row 156, shifting the existing historical rows (156-180) down to
# (157-181), matching the new weekly observation added upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 156:180 down to 157:181 by inserting a blank row at 156.
$ws.Rows.Item(156).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A156").Value = 7
$ws.Range("B156").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C156").Value = "Ñuble"
$ws.Range("D156").Value = 44474
$ws.Range("E156").Value = 16
$ws.Range("F156").Value = 100114013
$ws.Range("G156").Value = "Zanahoria"
$ws.Range("H156").Value = "Sin especificar"
$ws.Range("I156").Value = "Primera"
$ws.Range("J156").Value = 120
$ws.Range("K156").Value = 8000
$ws.Range("L156").Value = 9000
$ws.Range("M156").Value = 8500
$ws.Range("N156").Value = "$/saco 20 kilos"
$ws.Range("O156").Value = "Provincia de Diguillín"
$ws.Range("P156").Value = 425
$ws.Range("Q156").Value = 20
$ws.Range("R156").Value = "Hortaliza"
